# Iraq League workbook update (01-06-2024 01:16)
# The underlying source data had several match rows re-sorted; for a number
# of adjacent row pairs every column except "id" (column A) needs to be
# swapped between the two rows. We do this by swapping the B:AD range
# values (using Value2, which round-trips as a proper 2-D array through
# this COM host) between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")
    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2
    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

$rowPairs = @(
    @(17, 18),
    @(78, 79),
    @(89, 90),
    @(103, 104),
    @(173, 174),
    @(223, 224)
)

foreach ($pair in $rowPairs) {
    Swap-Rows $ws $pair[0] $pair[1]
}
